# Apply odds updates to Sheet1 as described by the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("T3").Value  = 9.5
$ws.Range("Y3").Value  = 35
$ws.Range("AH3").Value = 24
$ws.Range("AI3").Value = 18.5

# Row 11
$ws.Range("G11").Value  = 3.65
$ws.Range("H11").Value  = 3
$ws.Range("I11").Value  = 2.07
$ws.Range("L11").Value  = 1.45
$ws.Range("M11").Value  = 2.37
$ws.Range("Q11").Value  = 2.27
$ws.Range("R11").Value  = 2
$ws.Range("S11").Value  = 1.65
$ws.Range("Z11").Value  = 6.7
$ws.Range("AH11").Value = 19
$ws.Range("AI11").Value = 20

# Row 14
$ws.Range("G14").Value  = 3.55
$ws.Range("H14").Value  = 3.15
$ws.Range("I14").Value  = 2.05
$ws.Range("L14").Value  = 1.47
$ws.Range("M14").Value  = 2.32
$ws.Range("N14").Value  = 2.37
$ws.Range("O14").Value  = 1.45
$ws.Range("P14").Value  = 1.52
$ws.Range("Q14").Value  = 2.22
$ws.Range("T14").Value  = 7.8
$ws.Range("U14").Value  = 16.5
$ws.Range("Z14").Value  = 6.7
$ws.Range("AA14").Value = 6.3
$ws.Range("AB14").Value = 20
$ws.Range("AC14").Value = 150
$ws.Range("AF14").Value = 8.25
$ws.Range("AG14").Value = 9.5
$ws.Range("AH14").Value = 18
$ws.Range("AI14").Value = 21

# Row 16
$ws.Range("G16").Value  = 2.57
$ws.Range("K16").Value  = 8.5
$ws.Range("T16").Value  = 10.75
$ws.Range("Z16").Value  = 8.5
$ws.Range("AA16").Value = 7.1
$ws.Range("AE16").Value = 10.25
